$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.16"
$ws.Range("E2").Value = "'-2.40%"
$ws.Range("G2").Value = "'20"
$ws.Range("D3").Value = "'35.75"
$ws.Range("E3").Value = "'0.55%"
$ws.Range("G3").Value = "'20"
$ws.Range("D4").Value = "'5.083"
$ws.Range("E4").Value = "'-0.66%"
$ws.Range("G4").Value = "'20"
$ws.Range("D5").Value = "'0.08083"
$ws.Range("E5").Value = "'-1.55%"
$ws.Range("G5").Value = "'20"
$ws.Range("D6").Value = "'1.944"
$ws.Range("E6").Value = "'-6.16%"
$ws.Range("G6").Value = "'20"
$ws.Range("D7").Value = "'7.794"
$ws.Range("E7").Value = "'-2.49%"
$ws.Range("G7").Value = "'20"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9280"
$ws.Range("E8").Value = "'-0.01%"
$ws.Range("G8").Value = "'20"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1530"
$ws.Range("E9").Value = "'41.75%"
$ws.Range("G9").Value = "'20"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1896"
$ws.Range("E10").Value = "'-1.66%"
$ws.Range("G10").Value = "'20"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09004"
$ws.Range("E11").Value = "'-6.97%"
$ws.Range("G11").Value = "'20"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03443"
$ws.Range("E12").Value = "'-5.38%"
$ws.Range("G12").Value = "'20"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09923"
$ws.Range("E13").Value = "'-0.12%"
$ws.Range("G13").Value = "'20"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001427"
$ws.Range("E14").Value = "'-0.41%"
$ws.Range("G14").Value = "'20"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005761"
$ws.Range("E15").Value = "'1.49%"
$ws.Range("G15").Value = "'20"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.532"
$ws.Range("E16").Value = "'1.67%"
$ws.Range("G16").Value = "'20"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.050"
$ws.Range("E17").Value = "'-1.95%"
$ws.Range("G17").Value = "'20"
$ws.Range("D18").Value = "'2.990"
$ws.Range("E18").Value = "'0.31%"
$ws.Range("G18").Value = "'20"
$ws.Range("D19").Value = "'0.3442"
$ws.Range("E19").Value = "'0.46%"
$ws.Range("G19").Value = "'20"
$ws.Range("D20").Value = "'0.1299"
$ws.Range("E20").Value = "'-0.93%"
$ws.Range("G20").Value = "'20"
$ws.Range("D21").Value = "'5.036"
$ws.Range("E21").Value = "'-1.19%"
$ws.Range("G21").Value = "'20"
$ws.Range("D22").Value = "'0.2388"
$ws.Range("E22").Value = "'8.30%"
$ws.Range("G22").Value = "'20"
$ws.Range("D23").Value = "'0.04497"
$ws.Range("E23").Value = "'-0.99%"
$ws.Range("G23").Value = "'20"
$ws.Range("D24").Value = "'0.001208"
$ws.Range("E24").Value = "'-1.60%"
$ws.Range("G24").Value = "'20"
$ws.Range("D25").Value = "'0.004828"
$ws.Range("E25").Value = "'0.41%"
$ws.Range("G25").Value = "'20"
$ws.Range("D26").Value = "'0.0001225"
$ws.Range("E26").Value = "'-2.08%"
$ws.Range("G26").Value = "'20"
$ws.Range("D27").Value = "'0.0003015"
$ws.Range("E27").Value = "'-32.31%"
$ws.Range("G27").Value = "'20"
$ws.Range("G28").Value = "'20"
$ws.Range("G29").Value = "'20"
$ws.Range("G30").Value = "'20"
$ws.Range("G31").Value = "'20"
$ws.Range("G32").Value = "'20"
$ws.Range("G33").Value = "'20"
$ws.Range("G34").Value = "'20"
$ws.Range("G35").Value = "'20"
$ws.Range("G36").Value = "'20"
$ws.Range("G37").Value = "'20"
$ws.Range("G38").Value = "'20"
$ws.Range("D39").Value = "'0.01872"
$ws.Range("E39").Value = "'-6.17%"
$ws.Range("G39").Value = "'20"
$ws.Range("D40").Value = "'0.04810"
$ws.Range("E40").Value = "'-2.47%"
$ws.Range("G40").Value = "'20"
$ws.Range("D41").Value = "'0.01057"
$ws.Range("E41").Value = "'7.42%"
$ws.Range("G41").Value = "'20"
$ws.Range("D42").Value = "'0.007326"
$ws.Range("E42").Value = "'-4.44%"
$ws.Range("G42").Value = "'20"
$ws.Range("D43").Value = "'0.1349"
$ws.Range("E43").Value = "'-2.67%"
$ws.Range("G43").Value = "'20"
$ws.Range("D44").Value = "'0.002058"
$ws.Range("E44").Value = "'-2.80%"
$ws.Range("G44").Value = "'20"
$ws.Range("D45").Value = "'0.009718"
$ws.Range("E45").Value = "'-16.13%"
$ws.Range("G45").Value = "'20"
$ws.Range("D46").Value = "'0.00006224"
$ws.Range("E46").Value = "'-5.03%"
$ws.Range("G46").Value = "'20"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.40%"
$ws.Range("G47").Value = "'20"
$ws.Range("E48").Value = "'4.46%"
$ws.Range("G48").Value = "'20"
$ws.Range("G49").Value = "'20"
$ws.Range("D50").Value = "'0.00002094"
$ws.Range("E50").Value = "'-0.40%"
$ws.Range("G50").Value = "'20"
$ws.Range("D51").Value = "'0.0001994"
$ws.Range("E51").Value = "'-0.40%"
$ws.Range("G51").Value = "'20"
